$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.710.04'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.305.96'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.80'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.80'
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.506'
$ws.Range("E7").Value = '  +1.32%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  +0.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.49'
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("E11").Value = '  +6.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0784'
$ws.Range("E12").Value = '  +0.94%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.120'
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.72'
$ws.Range("E14").Value = '  +1.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.664.79'
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.308.88'
$ws.Range("E16").Value = '  +1.57%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.782'
$ws.Range("E17").Value = '  +2.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.632.01'
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.30'
$ws.Range("E19").Value = '  -2.37%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.06'
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0890'
$ws.Range("E21").Value = '  +0.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.95'
$ws.Range("E22").Value = '  +1.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.30'
$ws.Range("E23").Value = '  +5.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.76'
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("E26").Value = '  +1.01%  '
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.91'
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.10'
$ws.Range("E30").Value = '  +2.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.32'
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.01'
$ws.Range("E33").Value = '  +2.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.65'
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.46'
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0700'
$ws.Range("E36").Value = '  +3.31%  '
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0998'
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("E39").Value = '  +1.28%  '
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("E41").Value = '  +2.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.32'
$ws.Range("E42").Value = '  +15.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.954.86'
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.42'
$ws.Range("E44").Value = '  +4.49%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0278'
$ws.Range("E45").Value = '  +1.69%  '
$ws.Range("E46").Value = '  +2.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.74'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.533.86'
$ws.Range("E48").Value = '  +1.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.44'
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("E50").Value = '  -2.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.16'
$ws.Range("E51").Value = '  +2.60%  '
